$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without letting Excel
# auto-convert numeric-looking / percent-looking strings into numbers.
# We stage the text (forced via NumberFormat "@") in a scratch cell well
# outside the used range, copy it, and paste-special VALUES ONLY into the
# destination so no number format/style is transferred - matching the
# original inline-string cells which carry no explicit style.
function Set-TextValue($sheet, $cellRef, $text) {
    $staging = $sheet.Range("ZZ1")
    $staging.NumberFormat = "@"
    $staging.Value = $text
    $staging.Copy()
    $sheet.Range($cellRef).PasteSpecial(-4163, 0)
    $staging.Clear()
}

Set-TextValue $ws 'D2' '311.14'
Set-TextValue $ws 'E2' '1.10%'
Set-TextValue $ws 'D3' '37.48'
Set-TextValue $ws 'E3' '-0.46%'
Set-TextValue $ws 'D4' '5.111'
Set-TextValue $ws 'E4' '0.17%'
Set-TextValue $ws 'D5' '0.07836'
Set-TextValue $ws 'E5' '-0.42%'
Set-TextValue $ws 'D6' '1.917'
Set-TextValue $ws 'E6' '-3.62%'
Set-TextValue $ws 'D7' '8.226'
Set-TextValue $ws 'E7' '-0.17%'
Set-TextValue $ws 'D8' '2.725'
Set-TextValue $ws 'E8' '-12.96%'
Set-TextValue $ws 'D9' '0.9290'
Set-TextValue $ws 'E9' '0.45%'
Set-TextValue $ws 'D10' '0.1199'
Set-TextValue $ws 'E10' '-6.06%'
Set-TextValue $ws 'D11' '0.1900'
Set-TextValue $ws 'E11' '1.35%'
Set-TextValue $ws 'D12' '0.09320'
Set-TextValue $ws 'E12' '5.74%'
Set-TextValue $ws 'D13' '0.03425'
Set-TextValue $ws 'E13' '-0.10%'
Set-TextValue $ws 'D14' '0.09613'
Set-TextValue $ws 'E14' '-1.56%'
Set-TextValue $ws 'D15' '0.001365'
Set-TextValue $ws 'E15' '-1.93%'
Set-TextValue $ws 'D16' '0.005839'
Set-TextValue $ws 'E16' '-3.44%'
Set-TextValue $ws 'D17' '3.537'
Set-TextValue $ws 'E17' '-1.10%'
Set-TextValue $ws 'D18' '4.403'
Set-TextValue $ws 'E18' '1.33%'
Set-TextValue $ws 'E19' '-0.27%'
Set-TextValue $ws 'D20' '5.262'
Set-TextValue $ws 'E20' '5.02%'
Set-TextValue $ws 'D21' '0.1275'
Set-TextValue $ws 'E21' '-0.73%'
Set-TextValue $ws 'D22' '0.2587'
Set-TextValue $ws 'E22' '3.58%'
Set-TextValue $ws 'D23' '0.02104'
Set-TextValue $ws 'E23' '180.28%'
Set-TextValue $ws 'D24' '0.04351'
Set-TextValue $ws 'E24' '0.56%'
Set-TextValue $ws 'D25' '0.001196'
Set-TextValue $ws 'E25' '-2.13%'
Set-TextValue $ws 'D26' '0.004265'
Set-TextValue $ws 'E26' '-7.25%'
Set-TextValue $ws 'D27' '0.0001299'
Set-TextValue $ws 'E27' '-63.85%'
Set-TextValue $ws 'D39' '0.02085'
Set-TextValue $ws 'E39' '-9.46%'
Set-TextValue $ws 'D40' '0.05046'
Set-TextValue $ws 'E40' '0.33%'
Set-TextValue $ws 'D41' '0.007580'
Set-TextValue $ws 'E41' '0.67%'
Set-TextValue $ws 'D42' '0.009117'
Set-TextValue $ws 'E42' '-7.41%'
Set-TextValue $ws 'D43' '0.1349'
Set-TextValue $ws 'D44' '0.002003'
Set-TextValue $ws 'E44' '-4.33%'
Set-TextValue $ws 'D45' '0.008608'
Set-TextValue $ws 'E45' '7.07%'
Set-TextValue $ws 'D46' '0.00006699'
Set-TextValue $ws 'E46' '2.50%'
Set-TextValue $ws 'D47' '0.00000000749'
Set-TextValue $ws 'E47' '-0.40%'
Set-TextValue $ws 'B48' 'CoinbaseStockToken'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextValue $ws 'D48' '0.001199'
Set-TextValue $ws 'E48' '-0.41%'
Set-TextValue $ws 'B49' 'BOLO'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextValue $ws 'D49' '0.002910'
Set-TextValue $ws 'E49' '-3.19%'
Set-TextValue $ws 'D50' '0.00002098'
Set-TextValue $ws 'E50' '-0.40%'
Set-TextValue $ws 'D51' '0.0001998'
Set-TextValue $ws 'E51' '-0.40%'
